$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ticketDescriptionHighlighting annotations in column G (rows 2-4) were
# re-punctuated so the JSON-ish keys are properly quoted
# (start / end / key  ->  "start" / "end" / "key"). The logical highlighting
# data is unchanged, only the text formatting of the stored string.

$ws.Range("G2").Value2 = '[{" start": 143, "end": 147, "key": "System" }' + [char]0x200B + ',{ "start": 104, "end": 128, "key": "Fehlerbeschreibung" }' + [char]0x200B + ',{ "start": 67,"end": 77, "key": "System" }]'

$ws.Range("G3").Value2 = '[{"start": 229, "end": 297,"key": "Service Anfrage" },' + [char]0x200B + '{ "start": 191,"end": 192,"key": "System" }' + [char]0x200B + ',{"start": 176, "end": 191, "key": "System" }' + [char]0x200B + ',{"start": 129, "end": 144, "key": "System" }]'

$ws.Range("G4").Value2 = '[{"start": 130, "end": 165, "key": "Auslöser" }' + [char]0x200B + ',{ "start": 37, "end": 78, "key": "Fehlerbeschreibung" },' + [char]0x200B + '{ "start": 24, "end": 36, "key": "System" }]'

# The user had scrolled the sheet and moved the selection before saving
# (selection moved from F4 to G5).
$ws.Range("G5").Select()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
